$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Row 4
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1

# Row 5
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1

# Row 6
$ws.Range("H6").Value = 1

# Row 7
$ws.Range("H7").Value = 1

# Row 8
$ws.Range("H8").Value = 1

# Row 9
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 1

# Row 10
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1

# Row 11
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 1

# Row 12
$ws.Range("H12").Value = 1

# Row 13
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 1

# Row 14
$ws.Range("H14").Value = 1

# Row 15
$ws.Range("H15").Value = 1

# Row 16
$ws.Range("H16").Value = 1

# Row 17
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 1

# Row 18
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 1

$wb.Save()
